$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new column D ("unitType"), shifting the old D..M -> E..N.
#    The new column inherits the header-row style (border + fill) from its
#    neighbours automatically.
# ---------------------------------------------------------------------------
$ws.Columns("D:D").Insert()

# ---------------------------------------------------------------------------
# 2) Rename/rewrite the header (i18n) placeholders.
#    Before insert:  A date | B address | C unit | D accountNumber |
#                     E service | F number | G place | H value1 | I value2 |
#                     J value3 | K value4 | L clientName | M source
#    After insert:   C unit -> unitName, D (new) -> unitType,
#                     old "service" (now column F) -> resource
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "{d.i18n.unitName}"
$ws.Range("D1").Value = "{d.i18n.unitType}"
$ws.Range("F1").Value = "{d.i18n.resource}"

# ---------------------------------------------------------------------------
# 3) Add the new "unitType" placeholders to the two repeated meter rows.
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = "{d.meter[I].unitType}"
$ws.Range("D3").Value = "{d.meter[i + 1].unitType}"

# ---------------------------------------------------------------------------
# 4) Column widths (character-width units; Excel rounds to whole pixels, so
#    the inputs below are chosen to land on the closest achievable width to
#    the design spec).
# ---------------------------------------------------------------------------
$ws.Range("A1").ColumnWidth = 14.92
$ws.Range("B1").ColumnWidth = 19.25
$ws.Range("C1").ColumnWidth = 21.09
$ws.Range("D1").ColumnWidth = 21.09
$ws.Range("E1").ColumnWidth = 26.09
$ws.Range("F1").ColumnWidth = 20.09
$ws.Range("G1").ColumnWidth = 19.42
$ws.Range("H1").ColumnWidth = 17.25
$ws.Range("I1").ColumnWidth = 18.25
$ws.Range("J1").ColumnWidth = 18.25
$ws.Range("K1").ColumnWidth = 18.25
$ws.Range("L1").ColumnWidth = 18.25
$ws.Range("M1").ColumnWidth = 22.25
$ws.Range("N1").ColumnWidth = 18.42

# ---------------------------------------------------------------------------
# 5) Add a 7-row x 14-column sample/legend box (A4:N10) below the table:
#    white fill throughout, thin gray outline around the perimeter only.
# ---------------------------------------------------------------------------
$box = $ws.Range("A4:N10")
$box.Interior.Color = 16777215
$box.Interior.Pattern = 1
$box.BorderAround(1, 2, -4142, 11184810)

for ($r = 4; $r -le 10; $r++) {
    $ws.Rows($r).RowHeight = 15.35
}

Write-Host "edit applied"
